$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 29

# Set values in the same order the shared-string table gains new entries:
# D (Directed Graph) -> B (Graph) -> C (BFS, DFS)
$ws.Cells.Item($row, 1).Value = 28
$ws.Cells.Item($row, 4).Value = "Directed Graph"
$ws.Cells.Item($row, 2).Value = "Graph"
$ws.Cells.Item($row, 3).Value = "BFS, DFS"
$ws.Cells.Item($row, 5).Value = "easy"
$ws.Cells.Item($row, 6).Value = "GeeksForGeeks"

# Match formatting used by the other data rows: A/D/E/F centered, B left-aligned,
# C left at default (no explicit style), like the source row.
$ws.Cells.Item($row, 1).HorizontalAlignment = -4108
$ws.Cells.Item($row, 2).HorizontalAlignment = -4131
$ws.Cells.Item($row, 4).HorizontalAlignment = -4108
$ws.Cells.Item($row, 5).HorizontalAlignment = -4108
$ws.Cells.Item($row, 6).HorizontalAlignment = -4108

$ws.Range("C30").Select()
